# Update cryptocurrency price/volume figures per latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.706.14'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '2.099.82'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.26%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.388'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0841'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.48%  '
$ws.Range('D13').Value = '2.411.70'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.805'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '2.100.70'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '38.733.51'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '172.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('E27').Value = '  +1.73%  '
$ws.Range('E28').Value = '  +5.96%  '
$ws.Range('E29').Value = '  +4.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.97%  '
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0620'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.62%  '
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +4.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.22%  '
$ws.Range('D43').Value = '1.534.60'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.07%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0912'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').Value = '2.295.09'
$ws.Range('E51').Value = '  +0.09%  '
